$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the formatting of the existing (empty) trailing cell C20 before we
# overwrite row 20 with the new "language" entry, so it can be re-applied to
# the new trailing empty cell at C22.
$trailingFont = $ws.Range("C20").Font

# Add new "language" row at row 20 (A/B pick up the default column style,
# and C20 keeps the style it already had).
$ws.Range("A20").Value = "language"
$ws.Range("B20").Value = "Lenguaje"
$ws.Range("C20").Value = "Language"

# Re-create the trailing blank formatted cell, now moved down to row 22
# (row 21 is left unused), matching the formatting the old C20 cell had.
$ws.Range("C22").Font.Name = $trailingFont.Name
$ws.Range("C22").Font.Size = $trailingFont.Size
$ws.Range("C22").Font.Underline = $trailingFont.Underline
$ws.Range("C22").Font.Color = $trailingFont.Color
$ws.Range("C22").Font.Bold = $trailingFont.Bold
$ws.Range("C22").Font.Italic = $trailingFont.Italic

# Move the selection like in the diff (was C19, now B15)
$ws.Range("B15").Select()
